# Add data for 2021-10-03:
#  - Rename the sheet / header to reflect "through September 25" instead of "through September 24"
#  - Bump a handful of counts across several neighborhoods (rows) and months (columns)
#  - Add a few brand-new (previously empty) cells in column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename worksheet tab ---
$ws.Name = "Through 2021-09-25"

# --- Update the running-month column header (shared string) ---
$ws.Range("B1").Value = "September 2021 (through September 25)"

# --- Bump existing counts (row => column => new value) ---
$ws.Range("B2").Value = 15     # Garfield Park, Sep 2021 (through 25)
$ws.Range("K2").Value = 6      # Garfield Park, Sep 2020
$ws.Range("T2").Value = 5      # Garfield Park, Sep 2019
$ws.Range("AU2").Value = 2     # Garfield Park, Sep 2016

$ws.Range("B3").Value = 10     # North Lawndale, Sep 2021 (through 25)

$ws.Range("B5").Value = 10     # Austin, Sep 2021 (through 25)
$ws.Range("AL5").Value = 7     # Austin, Sep 2017

$ws.Range("AC6").Value = 5     # Roseland, Sep 2018

$ws.Range("AL7").Value = 2     # Auburn Gresham, Sep 2017

$ws.Range("AC8").Value = 2     # Little Village, Sep 2018

$ws.Range("B11").Value = 5     # Little Italy, UIC, Sep 2021 (through 25)
$ws.Range("K11").Value = 3     # Little Italy, UIC, Sep 2020

$ws.Range("B15").Value = 2     # Bucktown, Sep 2021 (through 25)

$ws.Range("B16").Value = 1     # West Pullman, Sep 2021 (through 25) -- new cell

$ws.Range("AL17").Value = 4    # South Shore, Sep 2017

$ws.Range("K23").Value = 6     # United Center, Sep 2020

$ws.Range("B25").Value = 1     # New City, Sep 2021 (through 25) -- new cell

$ws.Range("K36").Value = 2     # Calumet Heights, Sep 2020

$ws.Range("B45").Value = 1     # Albany Park, Sep 2021 (through 25) -- new cell

$ws.Range("K55").Value = 10    # Grand Crossing, Sep 2020

$ws.Range("B59").Value = 2     # Archer Heights, Sep 2021 (through 25)

$ws.Range("B65").Value = 3     # Chinatown, Sep 2021 (through 25)
